$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 125 (existing rows 125-138 shift down to 127-140).
$ws.Rows.Item(125).Insert()
$ws.Rows.Item(125).Insert()

# New row 125: Vega Monumental Concepcion - Naranja - New Hall - Primera
$ws.Cells.Item(125,1).Value  = 11
$ws.Cells.Item(125,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(125,3).Value  = "Bíobío"
$ws.Cells.Item(125,4).Value  = 44461
$ws.Cells.Item(125,5).Value  = 8
$ws.Cells.Item(125,6).Value  = "Fruta"
$ws.Cells.Item(125,7).Value  = 100102
$ws.Cells.Item(125,8).Value  = "Cítricos"
$ws.Cells.Item(125,9).Value  = 100102005
$ws.Cells.Item(125,10).Value = "Naranja"
$ws.Cells.Item(125,11).Value = "New Hall"
$ws.Cells.Item(125,12).Value = "Primera"
$ws.Cells.Item(125,13).Value = 100
$ws.Cells.Item(125,14).Value = 5500
$ws.Cells.Item(125,15).Value = 6000
$ws.Cells.Item(125,16).Value = 5750
$ws.Cells.Item(125,17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(125,18).Value = "Región de O'Higgins"
$ws.Cells.Item(125,19).Value = 383
$ws.Cells.Item(125,20).Value = 15

# New row 126: Vega Monumental Concepcion - Naranja - New Hall - Segunda
$ws.Cells.Item(126,1).Value  = 11
$ws.Cells.Item(126,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(126,3).Value  = "Bíobío"
$ws.Cells.Item(126,4).Value  = 44461
$ws.Cells.Item(126,5).Value  = 8
$ws.Cells.Item(126,6).Value  = "Fruta"
$ws.Cells.Item(126,7).Value  = 100102
$ws.Cells.Item(126,8).Value  = "Cítricos"
$ws.Cells.Item(126,9).Value  = 100102005
$ws.Cells.Item(126,10).Value = "Naranja"
$ws.Cells.Item(126,11).Value = "New Hall"
$ws.Cells.Item(126,12).Value = "Segunda"
$ws.Cells.Item(126,13).Value = 50
$ws.Cells.Item(126,14).Value = 4500
$ws.Cells.Item(126,15).Value = 4500
$ws.Cells.Item(126,16).Value = 4500
$ws.Cells.Item(126,17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(126,18).Value = "Región de O'Higgins"
$ws.Cells.Item(126,19).Value = 300
$ws.Cells.Item(126,20).Value = 15

Write-Output "edit complete"
